# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Newly computed "K" column values (column G) for data rows 2-57,
# replacing the previous Strike# derived values.
$sVals = @(1, 1, 3, 1, 1, 0, 0, 0, 1, 0, 3, 0, 1, 2, 0, 0, 1, 0, 1, 1, 2, 1, 1, 0, 2, 1, 2, 1, 0, 0, 1, 1, 1, 2, 0, 1, 0, 0, 0, 2, 0, 0, 1, 0, 2, 0, 1, 1, 1, 2, 2, 3, 2, 1, 1, 2)

$startRow = 2
for ($i = 0; $i -lt $sVals.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $sVals[$i]
}
